$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# Update the time_taken (F column) timestamps on the "data" sheet to reflect
# the re-run query time.
$ws1.Range("F2").Value = "2021-10-05 14:19:18.071010"
$ws1.Range("F3").Value = "2021-10-05 14:19:18.071018"
$ws1.Range("F4").Value = "2021-10-05 14:19:18.071021"
$ws1.Range("F5").Value = "2021-10-05 14:19:18.071023"
$ws1.Range("F6").Value = "2021-10-05 14:19:18.071026"
$ws1.Range("F7").Value = "2021-10-05 14:19:18.071029"
$ws1.Range("F8").Value = "2021-10-05 14:19:18.071032"
$ws1.Range("F9").Value = "2021-10-05 14:19:18.071034"
$ws1.Range("F10").Value = "2021-10-05 14:19:18.071037"
$ws1.Range("F11").Value = "2021-10-05 14:19:18.071040"
$ws1.Range("F12").Value = "2021-10-05 14:19:18.071042"
$ws1.Range("F13").Value = "2021-10-05 14:19:18.071045"
$ws1.Range("F14").Value = "2021-10-05 14:19:18.071047"
$ws1.Range("F15").Value = "2021-10-05 14:19:18.071049"
$ws1.Range("F16").Value = "2021-10-05 14:19:18.071052"
$ws1.Range("F17").Value = "2021-10-05 14:19:18.071054"
$ws1.Range("F18").Value = "2021-10-05 14:19:18.071057"
$ws1.Range("F19").Value = "2021-10-05 14:19:18.071060"
$ws1.Range("F20").Value = "2021-10-05 14:19:18.071062"
$ws1.Range("F21").Value = "2021-10-05 14:19:18.071064"
$ws1.Range("F22").Value = "2021-10-05 14:19:18.071067"
$ws1.Range("F23").Value = "2021-10-05 14:19:18.071069"
$ws1.Range("F24").Value = "2021-10-05 14:19:18.071072"
$ws1.Range("F25").Value = "2021-10-05 14:19:18.071074"
$ws1.Range("F26").Value = "2021-10-05 14:19:18.071077"
$ws1.Range("F27").Value = "2021-10-05 14:19:18.071080"
$ws1.Range("F28").Value = "2021-10-05 14:19:18.071082"
$ws1.Range("F29").Value = "2021-10-05 14:19:18.071085"
$ws1.Range("F30").Value = "2021-10-05 14:19:18.071087"

# Add a new "metadata" worksheet positioned after the "data" sheet.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Copy the header style (bold font + border + centered alignment) from the
# "data" sheet's header row, and the index-column style from its first data
# row, so the new sheet's formatting matches the rest of the workbook.
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row.
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row describing the source PanelApp query.
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Beckwith-Wiedemann syndrome (BWS) and other congenital overgrowth disorders"
$ws2.Range("C2").Value = 38
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.115"
$ws2.Range("E2").Value = "2021-08-11T12:15:32.977615Z"
$ws2.Range("F2").Value = "2021-10-05 14:19:18.067217"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/38/?format=json"

# Keep "data" as the active sheet/selection (unchanged by the source diff).
$ws1.Select() | Out-Null
$ws1.Range("A1").Select() | Out-Null
